# Efna5-Epha4 LR-pair sheet: refresh with updated TPM-derived values (rows 2-13, cols E-T)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2 = @(3, 1, 0.2708946666666667, 0.812684, 0.1616296696421007, 0.1616296696421007, 3, 1, 10.23061133333333, 30.691834, 0.4855635428718841, 0.4855635428718841, 2.771418046939556, 24.942762422456, 0.0784814750246306, 0.0784814750246306)
    3 = @(3, 1, 0.2708946666666667, 0.812684, 0.1616296696421007, 0.1616296696421007, 3, 1, 8.775186333333332, 26.325559, 0.4164864079521221, 0.4164864079521222, 2.377151176706222, 21.394360590356, 0.06731656052772667, 0.06731656052772667)
    4 = @(3, 1, 0.2708946666666667, 0.812684, 0.1616296696421007, 0.1616296696421007, 3, 1, 2.034752, 6.104255999999999, 0.09657305490303886, 0.09657305490303887, 0.5512034647893332, 4.960831183103999, 0.01560907096030662, 0.01560907096030662)
    5 = @(3, 1, 0.2708946666666667, 0.812684, 0.1616296696421007, 0.1616296696421007, 2, 0.6666666666666666, 0.02901266666666667, 0.087038, 0.001376994272954919, 0.001376994272954919, 0.007859376665777778, 0.070734389992, 0.0002225631294367681, 0.0002225631294367681)
    6 = @(2, 0.6666666666666666, 1.058121333333333, 3.174364, 0.6313295261673385, 0.6313295261673384, 3, 1, 10.23061133333333, 30.691834, 0.4855635428718841, 0.4855635428718841, 10.82522810484178, 97.427052943576, 0.3065506014454407, 0.3065506014454407)
    7 = @(2, 0.6666666666666666, 1.058121333333333, 3.174364, 0.6313295261673385, 0.6313295261673384, 3, 1, 8.775186333333332, 26.325559, 0.4164864079521221, 0.4164864079521222, 9.28521186327511, 83.56690676947601, 0.2629401665875501, 0.2629401665875501)
    8 = @(2, 0.6666666666666666, 1.058121333333333, 3.174364, 0.6313295261673385, 0.6313295261673384, 3, 1, 2.034752, 6.104255999999999, 0.09657305490303886, 0.09657305490303887, 2.153014499242667, 19.377130493184, 0.06096942099246789, 0.06096942099246789)
    9 = @(2, 0.6666666666666666, 1.058121333333333, 3.174364, 0.6313295261673385, 0.6313295261673384, 2, 0.6666666666666666, 0.02901266666666667, 0.087038, 0.001376994272954919, 0.001376994272954919, 0.03069892153688889, 0.276290293832, 0.0008693371418797675, 0.0008693371418797673)
    10 = @(2, 0.6666666666666666, 0.3470046666666667, 1.041014, 0.2070408041905609, 0.2070408041905609, 3, 1, 10.23061133333333, 30.691834, 0.4855635428718841, 0.4855635428718841, 3.550069875519556, 31.950628879676, 0.1005314664018128, 0.1005314664018128)
    11 = @(2, 0.6666666666666666, 0.3470046666666667, 1.041014, 0.2070408041905609, 0.2070408041905609, 3, 1, 8.775186333333332, 26.325559, 0.4164864079521221, 0.4164864079521222, 3.045030608536222, 27.405275476826, 0.08622968083684537, 0.08622968083684539)
    12 = @(2, 0.6666666666666666, 0.3470046666666667, 1.041014, 0.2070408041905609, 0.2070408041905609, 3, 1, 2.034752, 6.104255999999999, 0.09657305490303886, 0.09657305490303887, 0.7060684395093333, 6.354615955584, 0.01999456295026435, 0.01999456295026436)
    13 = @(2, 0.6666666666666666, 0.3470046666666667, 1.041014, 0.2070408041905609, 0.2070408041905609, 2, 0.6666666666666666, 0.02901266666666667, 0.087038, 0.001376994272954919, 0.001376994272954919, 0.01006753072577778, 0.09060777653200001, 0.0002850940016383831, 0.0002850940016383831)
}

foreach ($r in $newValues.Keys) {
    $vals = $newValues[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = 5 + $i  # column E = 5
        $ws.Cells.Item($r, $col).Value = $vals[$i]
    }
}

